# Revised Field Lesson Unit — apply commit-described edits via Word COM interop.
$d = $word.ActiveDocument

# ----------------------------------------------------------------------------
# 1) Title block: collapse "Assignment 08: ..." heading + "Total Points/Due
#    Date" paragraph into a single 4-line Heading1 paragraph.
# ----------------------------------------------------------------------------
$titleXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>Revised Field Lesson Unit</w:t><w:br/><w:t>TCE 486/586A</w:t><w:br/><w:t>Spring 2026 (Edwards)</w:t><w:br/><w:t>50 points possible</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertXML($titleXml)

# Re-apply bold to the title run only (excluding the trailing paragraph
# mark) so the serializer keeps an explicit <w:b/> instead of folding it
# into the (already-bold) Heading1 style.
$titlePara2 = $d.Paragraphs(1)
$titleBoldRng = $d.Range($titlePara2.Range.Start, $titlePara2.Range.End - 1)
$titleBoldRng.Font.Bold = 1

# Remove the old "Total Points: ... / Due Date: ..." paragraph that used to
# follow the heading; it is dropped entirely.
$d.Paragraphs(2).Range.Delete()

# ----------------------------------------------------------------------------
# 2) "Evaluation Criteria" heading becomes "3. How will my work be assessed?"
# ----------------------------------------------------------------------------
$evalHeadingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13) -eq "Evaluation Criteria") {
        $evalHeadingIndex = $i
        break
    }
}

$evalXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:rPr><w:color w:val="0F4761"/></w:rPr><w:t>3. How will my work be assessed?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$evalPara = $d.Paragraphs($evalHeadingIndex)
$evalPara.Range.InsertXML($evalXml)

$evalPara2 = $d.Paragraphs($evalHeadingIndex)
$evalBoldRng = $d.Range($evalPara2.Range.Start, $evalPara2.Range.End - 1)
$evalBoldRng.Font.Bold = 1

# ----------------------------------------------------------------------------
# 3) Replace everything from "Revised Lessons (15 points)" through the end
#    of "Final Thoughts" (i.e. the rest of the document body up to the
#    section break) with a single rubric table plus the new "4. Submission"
#    section. The bookmarkStart/bookmarkEnd markers that framed those
#    removed headings stay in place (now empty) because we only delete the
#    paragraph content between the heading we just rewrote and the very end
#    of the body range, leaving the sectPr untouched.
# ----------------------------------------------------------------------------
$bodyEnd = $d.Content.End - 1
$deleteStart = $d.Paragraphs($evalHeadingIndex + 1).Range.Start
$deleteRange = $d.Range($deleteStart, $bodyEnd)

$replacementXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:tbl><w:tblPr><w:tblW w:type="auto" w:w="0"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="3600"/><w:gridCol w:w="3600"/><w:gridCol w:w="3600"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Criterion</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Points</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>What I'm Looking For</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Completeness and polish</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>All lessons are thorough, well-organized, and professionally presented</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Inquiry depth</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Lessons demonstrate sophisticated understanding of inquiry teaching across multiple dimensions</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Integration of feedback</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Clear evidence that workshop feedback and peer review shaped meaningful revisions</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Technology showcase</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>One lesson demonstrates thoughtful tech integration that reveals or amplifies student thinking</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Differentiation showcase</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>One lesson demonstrates how you provide access to challenge for diverse learners</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Discourse showcase</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>One lesson demonstrates talk moves, questioning strategies, and structures that deepen student reasoning</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Unit coherence</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Lessons connect to each other and build toward meaningful learning</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Theoretical grounding</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Strong connections to course frameworks and readings with specific examples from your lessons</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Reflection depth</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="3600"/></w:tcPr><w:p><w:r><w:t>Thoughtful analysis of revision journey and growth as teacher in Unit Rationale</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:r><w:rPr><w:b/><w:color w:val="0F4761"/></w:rPr><w:br/><w:t>4. Submission</w:t></w:r></w:p><w:p><w:r><w:t>Submit a single PDF to Canvas containing:</w:t></w:r></w:p><w:p><w:r><w:t>&#8226; Cover page with your name and field placement context</w:t></w:r></w:p><w:p><w:r><w:t>&#8226; Unit rationale (4-5 pages)</w:t></w:r></w:p><w:p><w:r><w:t>&#8226; Three complete lesson plans</w:t></w:r></w:p><w:p><w:r><w:t>&#8226; Appendices (if needed): student handouts, assessments, supporting materials</w:t></w:r></w:p><w:p><w:r><w:br/><w:t>During Week 14, you'll also present a 10-13 minute Engagement Showcase from one of your revised lessons. This is evaluated separately as part of your participation grade.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$deleteRange.InsertXML($replacementXml)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count); Table count: $($d.Tables.Count)"
